$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# CaseField tab: insert a new "CategoryID" column (H) between
# FieldType (G) and FieldTypeParameter (old H, now I).
# ---------------------------------------------------------------------------
$caseField = $wb.Worksheets.Item("CaseField")
$caseField.Columns.Item(8).Insert()

$caseField.Range("H2").Value = "`t`nA non-mandatory field. This field will be used to indicate a category for a document field type or Collection of Document field ONLY. `nMust match to a valid CategoryID defined in the Categories tab for the given CaseTypeID."
$caseField.Range("H3").Value = "CategoryID"

# ---------------------------------------------------------------------------
# ComplexTypes tab: insert a new "CategoryID" column (F) between
# FieldType (E) and FieldTypeParameter (old F, now G).
# ---------------------------------------------------------------------------
$complexTypes = $wb.Worksheets.Item("ComplexTypes")
$complexTypes.Columns.Item(6).Insert()

$complexTypes.Range("F2").Value = "A non-mandatory field. This field will be used to indicate a category for a document field type or Collection of Document field ONLY. `nMust match to a valid CategoryID defined in the Categories tab."
$complexTypes.Range("F3").Value = "CategoryID"

# ---------------------------------------------------------------------------
# Restore cursor / view state roughly matching the authored workbook: the
# ComplexTypes tab was the active sheet, with the CategoryID column selected.
# ---------------------------------------------------------------------------
$caseField.Range("Q2").Select()
$complexTypes.Activate()
$complexTypes.Range("F3").Select()
